$wb = $excel.ActiveWorkbook

# Add the new "Washington County" worksheet after "Frederick County"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Washington County"

# Header row (mirrors Frederick County sheet layout)
$newSheet.Range("A1").Value = "Zone"
$newSheet.Range("B1").Value = "Zone Abbreviation"
$newSheet.Range("C1").Value = "Issuing Body"
$newSheet.Range("D1").Value = "Zone General Description"

# Environmental Conservation
$newSheet.Range("A2").Value = "Environmental Conservation"
$newSheet.Range("B2").Value = "EC"
$newSheet.Range("C2").Value = "Washington County Zoning Department"
$newSheet.Range("D2").Value = "The purpose of this district is to prescribe a zoning category for those areas where, because of natural geographic factors and existing land uses, it is considered feasible and desirable to conserve open spaces, water supply sources, woodland areas, wildlife and other natural resources. This district may include extensive steeply sloped areas, stream valleys, water supply sources, and wooded areas adjacent thereto."

# Preservation
$newSheet.Range("A3").Value = "Preservation"
$newSheet.Range("B3").Value = "P"
$newSheet.Range("C3").Value = "Washington County Zoning Department"
$newSheet.Range("D3").Value = "The purpose of this district is to prescribe a zoning category for those areas where, because of natural geographic factors and existing land uses, it is considered feasible and desirable to conserve open spaces, water supply sources, woodland areas, wildlife and other natural resources. This district includes the County’s designated Rural Legacy Area, federal lands, state parks, state wildlife management areas, county parks, Edgemont Watershed, and most of the mountaintops and the Potomac River."

# Column widths (best-fit, matches the XML diff)
$newSheet.Columns.Item(1).ColumnWidth = 26.28515625
$newSheet.Columns.Item(2).ColumnWidth = 17
$newSheet.Columns.Item(3).ColumnWidth = 36
$newSheet.Columns.Item(4).ColumnWidth = 23.5703125

# Move the new sheet after Frederick County
$freddy = $wb.Worksheets.Item("Frederick County")
$newSheet.Move($null, $freddy)

# Re-fetch references by name since Move() invalidates prior object handles
$freddy = $wb.Worksheets.Item("Frederick County")
$washington = $wb.Worksheets.Item("Washington County")

# Frederick County sheet view updates (zoom + selection), done while active
$freddy.Activate()
$excel.ActiveWindow.Zoom = 80
$freddy.Range("D27").Select()

# Activate the new Washington County sheet last so it is the active tab
$washington.Activate()
